$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has every (B,C) content pair shifted one row above
# its correct (A) label starting at row 13 (e.g. row13 A="Programa resumido:"
# but B/C hold the old "Semestral" placeholder, row18 A="Método:" but B/C
# hold the Docentes responsáveis name, etc). Fix this by inserting a new
# row at 13 - which pushes the existing (now correctly-labelled) rows 13-23
# down to 14-24 - and then filling in the newly revealed/blank cells with
# their real content.
$ws.Rows("13:13").Insert()

# Row 13 only has B/C content (no label in column A). The Insert() operation
# left a stray, empty A13 behind (it inherited formatting from row 12) -
# remove it completely so the row matches the target shape.
$ws.Range("A13").Clear()

# Give B13/C13 the normal body-text formatting used throughout the sheet
# (bold-label style for A, wrap-text style for B, red wrap-text style for C)
# by copying formats from the row above before writing the values.
$ws.Range("B11:C11").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13").Value = "230696 - Carlos José Todero Peixoto"
$ws.Range("C13").Value = "230696 - Carlos José Todero Peixoto"

# Objetivos / Objectives body text (row 10)
$ws.Range("B10").Value = "Verificação experimental das Leis da ótica e suas aplicações. Fenômenos físicos relativos à Física Moderna"
$ws.Range("C10").Value = "Verificação experimental das Leis da ótica e suas aplicações. Fenômenos físicos relativos à Física Moderna"

# Programa resumido (row 14, was row 13 before the insert)
$ws.Range("B14").Value = "Óptica geométrica e Física. Comprovações experimentais de física moderna."
$ws.Range("C14").Value = "Óptica geométrica e Física. Comprovações experimentais de física moderna."

# Programa (row 16, was row 15 before the insert)
$ws.Range("B16").Value = "1) Refração e reflexão.2) Espelhos planos e esféricos e lentes delgadas.3) Polarização.4) Interferência de ondas planas.5) Difração.6) Espectroscopia ótica.7) Determinação da constante de Planck.8) Radiação de corpo negro."
$ws.Range("C16").Value = "1) Refração e reflexão.2) Espelhos planos e esféricos e lentes delgadas.3) Polarização.4) Interferência de ondas planas.5) Difração.6) Espectroscopia ótica.7) Determinação da constante de Planck.8) Radiação de corpo negro."

# Método (row 19, was row 18 before the insert)
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Critério (row 20, was row 19 before the insert)
$ws.Range("B20").Value = "NF 5,0."
$ws.Range("C20").Value = "NF 5,0."

# Norma de recuperação (row 21, was row 20 before the insert)
$ws.Range("B21").Value = "(NF+RC)/2 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 5,0, onde RC é uma prova de recuperação a ser aplicada."

# Bibliografia (row 22, was row 21 before the insert)
$ws.Range("B22").Value = "Apostilas do Laboratório de Ensino de Física do IFSC/USP.RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 4, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 4, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 4, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 4, Thomson Pioneira (2008)."
$ws.Range("C22").Value = "Apostilas do Laboratório de Ensino de Física do IFSC/USP.RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 4, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 4, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 4, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 4, Thomson Pioneira (2008)."
